$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 605.8570999999999
$ws.Range("I33").Value = 226.57143
$ws.Range("J33").Value = 985.1429000000001
$ws.Range("K33").Value = 226.57143
$ws.Range("L33").Value = 985.1429000000001
$ws.Range("M33").Value = 2.428570000000008
$ws.Range("N33").Value = -1443.1429

$ws.Range("H38").Value = 818.94116
$ws.Range("I38").Value = 557.625
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 1672.875
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -1300.875
$ws.Range("N38").Value = -15744

$ws.Range("H98").Value = 1315.7142
$ws.Range("I98").Value = 1410
$ws.Range("J98").Value = 750
$ws.Range("K98").Value = 1410
$ws.Range("L98").Value = 750
$ws.Range("M98").Value = 88
$ws.Range("N98").Value = -3746

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1315.7142
$ws.Range("I122").Value = 1410
$ws.Range("J122").Value = 750
$ws.Range("K122").Value = 4230
$ws.Range("L122").Value = 2250
$ws.Range("M122").Value = -1780
$ws.Range("N122").Value = -7150

$ws.Range("H130").Value = 29164.834
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 29164.834
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 29164.834
$ws.Range("N130").Value = -39204.834

$ws.Range("H137").Value = 4468.0303
$ws.Range("I137").Value = 4381.567
$ws.Range("J137").Value = 5332.6665
$ws.Range("K137").Value = 13144.701
$ws.Range("L137").Value = 15997.9995
$ws.Range("M137").Value = -10594.701
$ws.Range("N137").Value = -21097.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13905.037
$ws.Range("I2").Value = 19136.834
$ws.Range("J2").Value = 3441.4443
$ws.Range("K2").Value = 19136.834
$ws.Range("L2").Value = 3441.4443
$ws.Range("M2").Value = -19023.834
$ws.Range("N2").Value = -3667.4443

$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -184

$ws.Range("H32").Value = 6031.689
$ws.Range("I32").Value = 6693.59
$ws.Range("J32").Value = 1729.3334
$ws.Range("K32").Value = 6693.59
$ws.Range("L32").Value = 1729.3334
$ws.Range("M32").Value = -6406.59
$ws.Range("N32").Value = -2303.3334

$ws.Range("H44").Value = 44747.25
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 44747.25
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 44747.25
$ws.Range("N44").Value = -45723.25

$ws.Range("H45").Value = 3697.5
$ws.Range("I45").Value = 3992.4285
$ws.Range("J45").Value = 3468.111
$ws.Range("K45").Value = 3992.4285
$ws.Range("L45").Value = 3468.111
$ws.Range("M45").Value = -3615.4285
$ws.Range("N45").Value = -4222.111

$ws.Range("H116").Value = 13905.037
$ws.Range("I116").Value = 19136.834
$ws.Range("J116").Value = 3441.4443
$ws.Range("K116").Value = 19136.834
$ws.Range("L116").Value = 3441.4443
$ws.Range("M116").Value = -16842.834
$ws.Range("N116").Value = -8029.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13905.037
$ws.Range("I3").Value = 19136.834
$ws.Range("J3").Value = 3441.4443
$ws.Range("K3").Value = 19136.834
$ws.Range("L3").Value = 3441.4443
$ws.Range("M3").Value = -19022.834
$ws.Range("N3").Value = -3669.4443

$ws.Range("H36").Value = 1099.2858
$ws.Range("I36").Value = 1099.2858
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1099.2858
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -565.2858000000001

$ws.Range("H99").Value = 103752.5
$ws.Range("I99").Value = 400010
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 400010
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -398512
$ws.Range("N99").Value = -7996

$ws.Range("H107").Value = 7238.4517
$ws.Range("I107").Value = 6295.16
$ws.Range("J107").Value = 11168.833
$ws.Range("K107").Value = 6295.16
$ws.Range("L107").Value = 11168.833
$ws.Range("M107").Value = -4375.16
$ws.Range("N107").Value = -15008.833

$ws.Range("H114").Value = 85310.5
$ws.Range("I114").Value = 70621
$ws.Range("J114").Value = 100000
$ws.Range("K114").Value = 70621
$ws.Range("L114").Value = 100000
$ws.Range("M114").Value = -66282
$ws.Range("N114").Value = -108678

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 322.5
$ws.Range("I7").Value = 157.5
$ws.Range("J7").Value = 432.5
$ws.Range("K7").Value = 157.5
$ws.Range("L7").Value = 432.5
$ws.Range("M7").Value = -44.5
$ws.Range("N7").Value = -658.5

$ws.Range("H16").Value = 1655.8
$ws.Range("I16").Value = 1548.5555
$ws.Range("J16").Value = 1816.6666
$ws.Range("K16").Value = 1548.5555
$ws.Range("L16").Value = 1816.6666
$ws.Range("M16").Value = -1261.5555
$ws.Range("N16").Value = -2390.6666

$ws.Range("H31").Value = 1607.0312
$ws.Range("I31").Value = 1607.0312
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1607.0312
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1312.0312

$ws.Range("H34").Value = 1607.0312
$ws.Range("I34").Value = 1607.0312
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1607.0312
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1405.0312

$ws.Range("H99").Value = 11813045
$ws.Range("I99").Value = 1745290.1
$ws.Range("J99").Value = 40002760
$ws.Range("K99").Value = 1745290.1
$ws.Range("L99").Value = 40002760
$ws.Range("M99").Value = -1743792.1
$ws.Range("N99").Value = -40005756

$ws.Range("H107").Value = 3769.4722
$ws.Range("I107").Value = 839.6
$ws.Range("J107").Value = 10428.272
$ws.Range("K107").Value = 839.6
$ws.Range("L107").Value = 10428.272
$ws.Range("M107").Value = 1080.4
$ws.Range("N107").Value = -14268.272

$ws.Range("H109").Value = 29666.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 29666.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 29666.5
$ws.Range("N109").Value = -31746.5

$ws.Range("H113").Value = 1655.8
$ws.Range("I113").Value = 1548.5555
$ws.Range("J113").Value = 1816.6666
$ws.Range("K113").Value = 1548.5555
$ws.Range("L113").Value = 1816.6666
$ws.Range("M113").Value = 621.4445000000001
$ws.Range("N113").Value = -6156.6666

$ws.Range("H126").Value = 11813045
$ws.Range("I126").Value = 1745290.1
$ws.Range("J126").Value = 40002760
$ws.Range("K126").Value = 5235870.300000001
$ws.Range("L126").Value = 120008280
$ws.Range("M126").Value = -5233400.300000001
$ws.Range("N126").Value = -120013220

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 10725
$ws.Range("I7").Value = 11422.223
$ws.Range("J7").Value = 4450
$ws.Range("K7").Value = 34266.669
$ws.Range("L7").Value = 13350
$ws.Range("M7").Value = -34154.669
$ws.Range("N7").Value = -13574

$ws.Range("H12").Value = 21599.6
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 26499.5
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 79498.5
$ws.Range("M12").Value = -5827
$ws.Range("N12").Value = -79844.5

$ws.Range("H86").Value = 1298.909
$ws.Range("I86").Value = 1737.8
$ws.Range("J86").Value = 933.1667
$ws.Range("K86").Value = 5213.4
$ws.Range("L86").Value = 2799.5001
$ws.Range("M86").Value = -4027.4
$ws.Range("N86").Value = -5171.5001

$ws.Range("H89").Value = 1298.909
$ws.Range("I89").Value = 1737.8
$ws.Range("J89").Value = 933.1667
$ws.Range("K89").Value = 15640.2
$ws.Range("L89").Value = 8398.5003
$ws.Range("M89").Value = -9712.199999999999
$ws.Range("N89").Value = -20254.5003

$ws.Range("H99").Value = 60849.832
$ws.Range("I99").Value = 2825
$ws.Range("J99").Value = 89862.25
$ws.Range("K99").Value = 8475
$ws.Range("L99").Value = 269586.75
$ws.Range("M99").Value = -6229
$ws.Range("N99").Value = -274078.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6238.981
$ws.Range("I132").Value = 5579.9346
$ws.Range("J132").Value = 10569.857
$ws.Range("K132").Value = 16739.8038
$ws.Range("L132").Value = 31709.571
$ws.Range("M132").Value = -14209.8038
$ws.Range("N132").Value = -36769.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11364368
$ws.Range("I22").Value = 30303598
$ws.Range("J22").Value = 829.6
$ws.Range("K22").Value = 30303598
$ws.Range("L22").Value = 829.6
$ws.Range("M22").Value = -30303303
$ws.Range("N22").Value = -1419.6

$ws.Range("H27").Value = 11364368
$ws.Range("I27").Value = 30303598
$ws.Range("J27").Value = 829.6
$ws.Range("K27").Value = 30303598
$ws.Range("L27").Value = 829.6
$ws.Range("M27").Value = -30303491
$ws.Range("N27").Value = -1043.6

$ws.Range("H33").Value = 9378.75
$ws.Range("I33").Value = 9007.5
$ws.Range("J33").Value = 9750
$ws.Range("K33").Value = 9007.5
$ws.Range("L33").Value = 9750
$ws.Range("M33").Value = -8717.5
$ws.Range("N33").Value = -10330

$ws.Range("H46").Value = 2852.1177
$ws.Range("I46").Value = 2449.8333
$ws.Range("J46").Value = 3071.5454
$ws.Range("K46").Value = 3071.5454
$ws.Range("L46").Value = 3071.5454
$ws.Range("M46").Value = -2261.8333
$ws.Range("N46").Value = -3447.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 59000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 59000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 59000
$ws.Range("N27").Value = -59138

$ws.Range("H62").Value = 3995.3333
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 4294.4
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 4294.4
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -5542.4

$ws.Range("H65").Value = 3995.3333
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 4294.4
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 21472
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -27712

$ws.Range("H107").Value = 376
$ws.Range("I107").Value = 394.41177
$ws.Range("J107").Value = 219.5
$ws.Range("K107").Value = 1183.23531
$ws.Range("L107").Value = 658.5
$ws.Range("M107").Value = 736.76469
$ws.Range("N107").Value = -4498.5

$ws.Range("H115").Value = 48428.145
$ws.Range("I115").Value = 49995
$ws.Range("J115").Value = 48307.617
$ws.Range("K115").Value = 49995
$ws.Range("L115").Value = 48307.617
$ws.Range("M115").Value = -48428
$ws.Range("N115").Value = -51441.617

$ws.Range("H136").Value = 1210.4615
$ws.Range("I136").Value = 1061.3334
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3184.0002
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -634.0001999999999
$ws.Range("N136").Value = -14100
